# Replace the "-" placeholder values in column C (sigma_ratio) with the
# numeric value 999 for the relevant rows, and move the active selection
# to I13 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(9, 10, 23, 28, 31)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 999
}

$ws.Range("I13").Select()
